$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column J (is_active) currently holds the formula =TRUE() (cached as
# numeric boolean 1) in every data row. The fix replaces that with the
# literal text value "TRUE" (a plain string, not a boolean/formula).
#
# Assigning the string "TRUE" straight to .Value/.Value2/.Formula gets
# auto-coerced back into a boolean by Excel's input parser, so instead we
# build the text result via a formula that evaluates to the *string*
# "TRUE" in a scratch cell, then copy/PasteSpecial-values it into every
# J cell - that preserves the literal string type (t="s") and leaves the
# cell's existing style untouched.
$scratch = $ws.Range("Z1")
$scratch.Formula = '=""&"TRUE"'
$scratch.Copy()

$firstRow = 2
$lastRow = 73
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 10).PasteSpecial(-4163)  # xlPasteValues
}

$scratch.ClearContents()

# Matches the diff's updated active-cell/selection (J1:J -> J2).
$ws.Range("J2").Select()
